$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.863.37'
$ws.Range("E2").Value = '  -2.14%  '

# Row 3
$ws.Range("D3").Value = '3.122.32'
$ws.Range("E3").Value = '  -0.17%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.26'
$ws.Range("E5").Value = '  -2.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.82'
$ws.Range("E6").Value = '  -4.92%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("D8").Value = '3.111.98'
$ws.Range("E8").Value = '  -0.42%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  -1.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  -4.01%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.23'
$ws.Range("E11").Value = '  -2.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.453'
$ws.Range("E12").Value = '  -3.14%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("E13").Value = '  -5.19%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.87'
$ws.Range("E14").Value = '  -3.63%  '

# Row 15
$ws.Range("D15").Value = '3.625.99'
$ws.Range("E15").Value = '  -0.37%  '

# Row 16
$ws.Range("E16").Value = '  +1.53%  '

# Row 17
$ws.Range("D17").Value = '62.833.19'
$ws.Range("E17").Value = '  -2.17%  '

# Row 18
$ws.Range("D18").Value = '3.116.13'
$ws.Range("E18").Value = '  -0.52%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.63'
$ws.Range("E19").Value = '  -3.74%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.38'
$ws.Range("E20").Value = '  -2.36%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.03'
$ws.Range("E21").Value = '  -3.58%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.693'
$ws.Range("E22").Value = '  -2.55%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.63'
$ws.Range("E23").Value = '  -0.65%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.25'
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.85'
$ws.Range("E25").Value = '  -4.04%  '

# Row 26
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("E27").Value = '  -1.91%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.82'
$ws.Range("E28").Value = '  -6.36%  '

# Row 29
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.08'
$ws.Range("E29").Value = '  +1.64%  '

# Row 30
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.81'
$ws.Range("E30").Value = '  -4.74%  '

# Row 31
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.53'
$ws.Range("E32").Value = '  -1.29%  '

# Row 33
$ws.Range("E33").Value = '  -5.59%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.51'
$ws.Range("E34").Value = '  -4.71%  '

# Row 35
$ws.Range("E35").Value = '  -3.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.87'
$ws.Range("E36").Value = '  -0.92%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.72'
$ws.Range("E37").Value = '  -3.95%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0675'
$ws.Range("E38").Value = '  -12.11%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0384'
$ws.Range("E39").Value = '  -2.20%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '415.74'
$ws.Range("E40").Value = '  -6.65%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.15'
$ws.Range("E41").Value = '  -0.51%  '

# Row 42
$ws.Range("D42").Value = '2.889.96'
$ws.Range("E42").Value = '  +1.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.65'
$ws.Range("E43").Value = '  -11.93%  '

# Row 44
$ws.Range("E44").Value = '  -6.38%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.09'
$ws.Range("E47").Value = '  -5.80%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.27'
$ws.Range("E48").Value = '  -2.78%  '

# Row 49
$ws.Range("E49").Value = '  -0.75%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.23'
$ws.Range("E50").Value = '  -8.29%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.27'
$ws.Range("E51").Value = '  -0.02%  '
